# Daily attendance processing - 2025-10-30 17:20:38
#
# Normalizes the "Recorded By" column (column G) on the active sheet.
# Whenever the literal token "System" (capital S) appears in the
# comma-separated list of recorders but is not already the first
# entry, it is moved to the front of the list:
#   - if a lowercase "system" token also exists in the same list, the
#     two tokens simply swap places (preserving both, with "System"
#     now first and "system" taking the old "System" slot)
#   - otherwise "System" is relocated to the front and the remaining
#     entries keep their original relative order
#
# Rows that are blank, already start with "System", or do not contain
# an exact "System" token are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($val -eq $null -or $val -eq "") {
        continue
    }

    $rawParts = $val.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    # Locate exact-case tokens "system" and "System" (ordinal/case-sensitive
    # compare via CompareTo, since -eq/-ceq are case-insensitive here).
    $lowerIdx = -1
    $upperIdx = -1
    for ($i = 0; $i -lt $parts.Count; $i++) {
        if ($parts[$i].CompareTo("system") -eq 0) {
            $lowerIdx = $i
        }
        if ($parts[$i].CompareTo("System") -eq 0) {
            $upperIdx = $i
        }
    }

    $changed = $false
    $newVal = $val

    if ($lowerIdx -ge 0 -and $upperIdx -ge 0) {
        # Both "system" and "System" present -> swap their positions.
        $tmp = $parts[$lowerIdx]
        $parts[$lowerIdx] = $parts[$upperIdx]
        $parts[$upperIdx] = $tmp
        $newVal = $parts -join ", "
        $changed = $true
    } elseif ($upperIdx -gt 0) {
        # "System" present but not first -> move it to the front.
        $sysTok = $parts[$upperIdx]
        $newParts = @()
        for ($i = 0; $i -lt $parts.Count; $i++) {
            if ($i -ne $upperIdx) {
                $newParts += $parts[$i]
            }
        }
        $finalParts = @($sysTok) + $newParts
        $newVal = $finalParts -join ", "
        $changed = $true
    }

    if ($changed) {
        $cell.Value = $newVal
    }
}
